# Adds " React Native," after "React.js," in three places in the résumé,
# splitting runs so the new text lands in its own run (matching the
# original run-per-phrase structure of the document).

$d = $word.ActiveDocument

function Split-AfterInsert($paraRange, [string[]]$segments) {
    # Re-establishes run boundaries along $segments (in left-to-right
    # order) inside $paraRange, which the preceding Insert call may have
    # coalesced into one big run. Toggling Bold on/off for each exact
    # sub-range forces that sub-range into its own run without leaving any
    # residual formatting behind (the net Bold change is a no-op).
    foreach ($seg in $segments) {
        $probe = $paraRange.Duplicate
        $ok = $probe.Find.Execute($seg, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
        if ($ok) {
            $probe.Bold = 1
            $probe.Bold = 0
        }
    }
}

# ---------------------------------------------------------------------
# Edit 1: top summary line — plain runs, no rPr.
# ---------------------------------------------------------------------
$para1 = $d.Paragraphs(7).Range
$f1 = $para1.Duplicate
$f1.Find.Execute("React.js,", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$ins1 = $d.Range($f1.End, $f1.End)
$ins1.InsertAfter(" React Native,")

Split-AfterInsert $d.Paragraphs(7).Range @(
    "Experienced Full Stack Developer | 3+ Years Expertise in React.js,",
    " React Native,",
    " Next.js, Node.js, TypeScript, MongoDB,",
    " MySQL, Postgres,",
    " API Integration, and Web App Development"
)

# ---------------------------------------------------------------------
# Edit 2: professional-summary paragraph — Calibri + white-shaded runs.
# ---------------------------------------------------------------------
$para2 = $d.Paragraphs(11).Range
$f2 = $para2.Duplicate
$f2.Find.Execute("React.js,", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$ins2 = $d.Range($f2.End, $f2.End)
$ins2.InsertAfter(" React Native,")

Split-AfterInsert $d.Paragraphs(11).Range @(
    "React.js,",
    " React Native,",
    " Next"
)

# ---------------------------------------------------------------------
# Edit 3: skills line — Calibri runs (no shading).
# ---------------------------------------------------------------------
$para3 = $d.Paragraphs(21).Range
$f3 = $para3.Duplicate
$f3.Find.Execute("React.JS, T", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$insertAt = $f3.End - 1
$ins3 = $d.Range($insertAt, $insertAt)
$ins3.InsertBefore("React Native, ")

Split-AfterInsert $d.Paragraphs(21).Range @(
    "JavaScript,",
    " React.JS, ",
    "React Native, ",
    "Typescript, Next.JS,",
    " Node.js, Express.js, Strapi,",
    " ECMASCRIPT 5",
    "/ES6."
)

Write-Host "Done"
